$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 B11: change style so it matches the "normal" date format (style index 2)
$ws.Range("B11").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Add new row 12
$ws.Range("A12").Value = 806.651
$ws.Range("B12").Value = 45742
$ws.Range("B12").NumberFormat = "YYYY-MM-DD"
$ws.Range("C12").Value = 773.9160000000001
